# daily auto push: 2026-01-14 22:36 UTC
#
# The tracker sheet (excel/sei2.xlsx -> Sheet1) logs one row per
# (date, weekday, hour, ranking) observation, in chronological order.
# A new observation for 2026/01/15 at hour 6 arrived and needs to be
# inserted right after the existing "2026/01/15" rows (which end at
# row 640), pushing every following row down by one. The sheet's used
# range grows from A1:D682 to A1:D683.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 641; everything from the old row 641
# onward (2026/12/29 ... 2027/01/05) shifts down to 642..683.
$ws.Rows("641:641").Insert()

# Column A holds dates formatted as plain text (e.g. "2026/01/15"),
# matching every other row in the column. Force text formatting before
# assigning so the engine doesn't auto-coerce the string into a date
# serial number, then drop the temporary format override so the cell
# is left with the same (default) style as its neighbours.
$ws.Range("A641").NumberFormat = "@"
$ws.Range("A641").Value = "2026/01/15"
$ws.Range("A641").ClearFormats()

$ws.Range("B641").Value = "木"
$ws.Range("C641").Value = 6
$ws.Range("D641").Value = 201
